$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Sheet2")

# Copy Sheet1's row 3 (A12,B12,C12,D12) and what was row 4 (A13,B13,C13,D13)
# down onto Sheet2 as new rows 6 and 7 (row 5 stays blank).
$ws2.Range("A6").Value = "A 12"
$ws2.Range("B6").Value = "B 12"
$ws2.Range("C6").Value = "C 12"
$ws2.Range("D6").Value = "D 12"
$ws2.Rows.Item(6).RowHeight = 13.8

$ws2.Range("A7").Value = "A 13"
$ws2.Range("B7").Value = "B 13"
$ws2.Range("C7").Value = "C 13"
$ws2.Range("D7").Value = "D 13"
$ws2.Rows.Item(7).RowHeight = 13.8

# Remove the now-duplicated last row from Sheet1.
$ws1.Rows.Item(4).Delete()

# Tweak Sheet1's row3 height / column C width slightly (cosmetic).
$ws1.Rows.Item(3).RowHeight = 13.8
$ws1.Columns.Item(3).ColumnWidth = 14.2

# Sheet1 becomes the active sheet/tab; move the on-sheet selections too.
$null = $ws2.Range("A6").Select()
$ws2.Select()
$null = $ws1.Range("A7").Select()
$ws1.Activate()
